$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each entry: target cell reference and the new text value.
# NumberFormat is forced to Text ("@") before assignment so that
# numeric-looking strings (prices, percentages, hour values) are
# stored as text, matching the original inline-string content type
# instead of being auto-converted to numbers/percentages by Excel.
$changes = @(
    @{ Cell = "D2"; Value = '272.98' },
    @{ Cell = "E2"; Value = '1.23%' },
    @{ Cell = "G2"; Value = '21' },
    @{ Cell = "D3"; Value = '26.80' },
    @{ Cell = "E3"; Value = '0.38%' },
    @{ Cell = "G3"; Value = '21' },
    @{ Cell = "D4"; Value = '4.913' },
    @{ Cell = "E4"; Value = '4.36%' },
    @{ Cell = "G4"; Value = '21' },
    @{ Cell = "D5"; Value = '0.06325' },
    @{ Cell = "E5"; Value = '3.66%' },
    @{ Cell = "G5"; Value = '21' },
    @{ Cell = "D6"; Value = '6.942' },
    @{ Cell = "E6"; Value = '3.03%' },
    @{ Cell = "G6"; Value = '21' },
    @{ Cell = "D7"; Value = '3.355' },
    @{ Cell = "E7"; Value = '5.93%' },
    @{ Cell = "G7"; Value = '21' },
    @{ Cell = "D8"; Value = '1.369' },
    @{ Cell = "E8"; Value = '53.61%' },
    @{ Cell = "G8"; Value = '21' },
    @{ Cell = "D9"; Value = '0.8870' },
    @{ Cell = "E9"; Value = '3.38%' },
    @{ Cell = "G9"; Value = '21' },
    @{ Cell = "D10"; Value = '0.1474' },
    @{ Cell = "E10"; Value = '3.74%' },
    @{ Cell = "G10"; Value = '21' },
    @{ Cell = "D11"; Value = '0.05133' },
    @{ Cell = "E11"; Value = '2.11%' },
    @{ Cell = "G11"; Value = '21' },
    @{ Cell = "D12"; Value = '0.07337' },
    @{ Cell = "E12"; Value = '3.26%' },
    @{ Cell = "G12"; Value = '21' },
    @{ Cell = "D13"; Value = '0.03169' },
    @{ Cell = "E13"; Value = '-0.14%' },
    @{ Cell = "G13"; Value = '21' },
    @{ Cell = "D14"; Value = '0.09050' },
    @{ Cell = "E14"; Value = '0.21%' },
    @{ Cell = "G14"; Value = '21' },
    @{ Cell = "D15"; Value = '0.001579' },
    @{ Cell = "E15"; Value = '2.92%' },
    @{ Cell = "G15"; Value = '21' },
    @{ Cell = "D16"; Value = '0.0006346' },
    @{ Cell = "E16"; Value = '4.73%' },
    @{ Cell = "G16"; Value = '21' },
    @{ Cell = "D17"; Value = '0.006024' },
    @{ Cell = "E17"; Value = '-0.34%' },
    @{ Cell = "G17"; Value = '21' },
    @{ Cell = "D18"; Value = '3.478' },
    @{ Cell = "E18"; Value = '0.43%' },
    @{ Cell = "G18"; Value = '21' },
    @{ Cell = "E19"; Value = '1.65%' },
    @{ Cell = "G19"; Value = '21' },
    @{ Cell = "E20"; Value = '2.22%' },
    @{ Cell = "G20"; Value = '21' },
    @{ Cell = "D21"; Value = '0.1335' },
    @{ Cell = "E21"; Value = '2.81%' },
    @{ Cell = "G21"; Value = '21' },
    @{ Cell = "D22"; Value = '3.915' },
    @{ Cell = "E22"; Value = '1.98%' },
    @{ Cell = "G22"; Value = '21' },
    @{ Cell = "D23"; Value = '0.04338' },
    @{ Cell = "E23"; Value = '2.22%' },
    @{ Cell = "G23"; Value = '21' },
    @{ Cell = "D24"; Value = '0.001182' },
    @{ Cell = "E24"; Value = '-0.32%' },
    @{ Cell = "G24"; Value = '21' },
    @{ Cell = "D25"; Value = '0.003645' },
    @{ Cell = "E25"; Value = '-12.15%' },
    @{ Cell = "G25"; Value = '21' },
    @{ Cell = "D26"; Value = '0.0001203' },
    @{ Cell = "E26"; Value = '0.20%' },
    @{ Cell = "G26"; Value = '21' },
    @{ Cell = "E27"; Value = '15.53%' },
    @{ Cell = "G27"; Value = '21' },
    @{ Cell = "G28"; Value = '21' },
    @{ Cell = "G29"; Value = '21' },
    @{ Cell = "G30"; Value = '21' },
    @{ Cell = "G31"; Value = '21' },
    @{ Cell = "G32"; Value = '21' },
    @{ Cell = "G33"; Value = '21' },
    @{ Cell = "G34"; Value = '21' },
    @{ Cell = "G35"; Value = '21' },
    @{ Cell = "G36"; Value = '21' },
    @{ Cell = "G37"; Value = '21' },
    @{ Cell = "G38"; Value = '21' },
    @{ Cell = "G39"; Value = '21' },
    @{ Cell = "D40"; Value = '0.04026' },
    @{ Cell = "E40"; Value = '1.97%' },
    @{ Cell = "G40"; Value = '21' },
    @{ Cell = "D41"; Value = '0.006621' },
    @{ Cell = "E41"; Value = '58.40%' },
    @{ Cell = "G41"; Value = '21' },
    @{ Cell = "D42"; Value = '0.1163' },
    @{ Cell = "E42"; Value = '4.09%' },
    @{ Cell = "G42"; Value = '21' },
    @{ Cell = "D43"; Value = '0.002355' },
    @{ Cell = "E43"; Value = '17.17%' },
    @{ Cell = "G43"; Value = '21' },
    @{ Cell = "D44"; Value = '0.01258' },
    @{ Cell = "E44"; Value = '-1.09%' },
    @{ Cell = "G44"; Value = '21' },
    @{ Cell = "D45"; Value = '0.00005270' },
    @{ Cell = "E45"; Value = '2.97%' },
    @{ Cell = "G45"; Value = '21' },
    @{ Cell = "E46"; Value = '132.72%' },
    @{ Cell = "G46"; Value = '21' },
    @{ Cell = "E47"; Value = '-13.20%' },
    @{ Cell = "G47"; Value = '21' },
    @{ Cell = "E48"; Value = '-0.03%' },
    @{ Cell = "G48"; Value = '21' },
    @{ Cell = "G49"; Value = '21' },
    @{ Cell = "G50"; Value = '21' },
    @{ Cell = "G51"; Value = '21' }
)

foreach ($chg in $changes) {
    $rng = $ws.Range($chg.Cell)
    $rng.NumberFormat = "@"
    $rng.Value = $chg.Value
}
